$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindContinue = 1
$wdReplaceAll = 2
$wdFindContinue = 1

# "find /usr/bin -name *td"   -> "find /usr/bin *td"
# "find /usr/bin -name *td*"  -> "find /usr/bin *td*"
# "find /usr/bin -name *td || find /usr/bin -name *cd" ->
#   "find /usr/bin *td || find /usr/bin *cd"   (first pass)
$r1 = $d.Content
$r1.Find.Execute("-name ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceAll)

# "find /usr/bin *td || find /usr/bin *cd" -> "find /usr/bin *td && find /usr/bin *cd"
$r2 = $d.Content
$r2.Find.Execute("||", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "&&", $wdReplaceAll)
